# update al 02 de octubre 2023
$wb = $excel.ActiveWorkbook

$ingreso = $wb.Worksheets.Item("Ingreso")
$gastos  = $wb.Worksheets.Item("Gastos")
$cxc     = $wb.Worksheets.Item("Cuentas por cobrar")

# --- Ingreso: new "Aporte" rows for 2023-10-01 (serial 45200) ---
$ingresoRows = @(
    @(45200, "Carlos",  300,  "Aporte"),
    @(45200, "Javier",  1200, "Aporte"),
    @(45200, "Johan",   600,  "Aporte"),
    @(45200, "Julio",   100,  "Aporte"),
    @(45200, "Invitados", 100, "Aporte"),
    @(45200, "chamo",   200,  "Aporte"),
    @(45200, "Kawai",   50,   "Aporte"),
    @(45200, "Joel",    300,  "Aporte"),
    @(45200, "Frandy",  1000, "Aporte")
)

$startRow = 533
for ($i = 0; $i -lt $ingresoRows.Count; $i++) {
    $r = $startRow + $i
    $row = $ingresoRows[$i]
    $ingreso.Range("A$r").Value = $row[0]
    $ingreso.Range("B$r").Value = $row[1]
    $ingreso.Range("C$r").Value = $row[2]
    $ingreso.Range("C$r").Style = "Normal"
    $ingreso.Range("D$r").Value = $row[3]
}

# --- Gastos: new "Arbitro, agua y hielo" rows ---
$gastos.Range("A61").Value = 45193
$gastos.Range("B61").Value = "Arbitro, agua y hielo"
$gastos.Range("C61").Value = 950

$gastos.Range("A62").Value = 45200
$gastos.Range("B62").Value = "Arbitro, agua y hielo"
$gastos.Range("C62").Value = 940

# --- view / selection state ---
# Gastos selection moves to C62 (sheet stops being the active tab)
[void]$gastos.Range("C62").Select()

# "Cuentas por cobrar" selection moves to D10
[void]$cxc.Range("D10").Select()

# Ingreso becomes the active sheet/tab again, with final selection on B534
[void]$ingreso.Range("B534").Select()
